# Update "想去人数" (F column) values across the four worksheets to reflect
# the latest generated output (commit: "Update gh-pages to output generated
# at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览" = @(
        @(2, 4529), @(5, 3632), @(10, 355), @(11, 2511), @(13, 36),
        @(14, 1967), @(15, 274), @(16, 15), @(17, 552), @(20, 10381),
        @(21, 6044), @(22, 17), @(23, 8), @(26, 5), @(27, 10),
        @(28, 838), @(29, 18), @(30, 170), @(31, 857), @(35, 476),
        @(36, 123), @(37, 263), @(39, 242), @(40, 4845), @(41, 25),
        @(42, 1130), @(43, 165), @(44, 174), @(45, 97), @(46, 485)
    )
    "演出" = @(
        @(15, 3559), @(25, 37)
    )
    "本地生活" = @(
        @(2, 8792), @(3, 440), @(4, 1632)
    )
    "全部类型" = @(
        @(2, 440), @(3, 1632), @(5, 4529), @(8, 3632), @(12, 355),
        @(13, 2511), @(17, 36), @(18, 274), @(19, 15), @(21, 552),
        @(24, 10381), @(25, 3559), @(27, 17), @(30, 5), @(31, 10),
        @(32, 838), @(33, 18), @(34, 170), @(35, 857), @(38, 123),
        @(39, 263), @(41, 242), @(42, 4845), @(43, 25), @(44, 1130),
        @(45, 165), @(46, 97), @(47, 485)
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $updates[$sheetName]) {
        $row = $pair[0]
        $newValue = $pair[1]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}
